# Applies the "cryptos list" update described by the commit diff.
# The sheet stores Price (col D) and Volume 1h (col E) as plain text cells
# (t="inlineStr" in the original OOXML), never as real numbers/percentages.
# We must replicate that: most new values are not number-like (contain two
# dots, a trailing "%", or surrounding spaces) so a normal .Value assignment
# is stored as text automatically. A handful of new Price values DO look like
# plain numbers (e.g. "1.0000", "0.07200", "306.11") - assigning those naively
# would make Excel coerce them into numeric cells and silently drop trailing
# zeros / precision. For those we force text by prefixing the value with a
# leading apostrophe (the standard "treat as text" marker) and then reset the
# cell Style back to "Normal" so no stray number-format style sticks around.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.Value = $text
}

function Set-ForcedTextValue($range, $text) {
    # Leading apostrophe forces Excel to store the value as text even though
    # it looks like a number; resetting the style afterwards avoids leaving
    # a "text number format" behind on the cell.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

$textUpdates = @{
    "D2" = '27.026.89'
    "E2" = '  +0.48%  '
    "D3" = '1.889.46'
    "E3" = '  +1.45%  '
    "E4" = '  +0.08%  '
    "E5" = '  +0.43%  '
    "E6" = '  +0.08%  '
    "E7" = '  +2.80%  '
    "E8" = '  +2.86%  '
    "E9" = '  +0.39%  '
    "E10" = '  +2.08%  '
    "E11" = '  +0.89%  '
    "E12" = '  +1.81%  '
    "D13" = '1.887.53'
    "E13" = '  +1.19%  '
    "E14" = '  -0.69%  '
    "E15" = '  +0.09%  '
    "E16" = '  +0.10%  '
    "E17" = '  -0.08%  '
    "E18" = '  +1.67%  '
    "E19" = '  -0.02%  '
    "D20" = '27.070.30'
    "E20" = '  +0.45%  '
    "E21" = '  +0.47%  '
    "D22" = '2.116.17'
    "E22" = '  -0.33%  '
    "E23" = '  +1.94%  '
    "E24" = '  -0.30%  '
    "E25" = '  +9.91%  '
    "E26" = '  -1.33%  '
    "E27" = '  -2.52%  '
    "E28" = '  +0.98%  '
    "E29" = '  +0.75%  '
    "E30" = '  +5.23%  '
    "E31" = '  +1.78%  '
    "E32" = '  -0.40%  '
    "E33" = '  -2.03%  '
    "E34" = '  +7.41%  '
    "E35" = '  +2.16%  '
    "E36" = '  -0.02%  '
    "E37" = '  +0.51%  '
    "E38" = '  +0.45%  '
    "E39" = '  +0.44%  '
    "E40" = '  -0.67%  '
    "E41" = '  +0.44%  '
    "E42" = '  +5.03%  '
    "E43" = '  +0.66%  '
    "E44" = '  +2.06%  '
    "E45" = '  +2.20%  '
    "E46" = '  +2.75%  '
    "B47" = 'EnergySwap'
    "C47" = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    "E47" = '  +1.59%  '
    "B48" = 'PaxDollar'
    "C48" = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
    "E48" = '  +0.11%  '
    "E49" = '  +2.23%  '
    "E50" = '  +2.66%  '
    "E51" = '  +1.40%  '
}

$numericLookingUpdates = @{
    "D5" = '306.11'
    "D7" = '0.5189'
    "D8" = '0.3753'
    "D9" = '0.07200'
    "D10" = '21.12'
    "D11" = '0.9012'
    "D12" = '0.07634'
    "D14" = '94.48'
    "D15" = '5.239'
    "D16" = '1.002'
    "D17" = '0.000008503'
    "D19" = '1.0000'
    "D21" = '5.050'
    "D24" = '6.389'
    "D26" = '145.81'
    "D27" = '1.736'
    "D28" = '18.07'
    "D29" = '114.20'
    "D30" = '4.914'
    "D31" = '4.788'
    "D32" = '0.09188'
    "D33" = '0.05040'
    "D34" = '1.240'
    "D35" = '0.7668'
    "D36" = '2.960'
    "D37" = '3.274'
    "D38" = '2.604'
    "D39" = '0.5598'
    "D40" = '0.01988'
    "D41" = '1.074'
    "D42" = '9.012'
    "D43" = '6.618'
    "D44" = '118.76'
    "D46" = '0.4826'
    "D47" = '10.19'
    "D48" = '1.000'
    "D50" = '37.68'
}

foreach ($ref in $textUpdates.Keys) {
    Set-TextValue $ws.Range($ref) $textUpdates[$ref]
}

foreach ($ref in $numericLookingUpdates.Keys) {
    Set-ForcedTextValue $ws.Range($ref) $numericLookingUpdates[$ref]
}
